$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(42602.576018518521, "Named", 3115, 408, 16, 6, 7, 46, 53, 0, 1, 0, 100),
    @(42602.576863425929, "Named", 3118, 408, 16, 6, 7, 46, 53, 0, 1, 0, 100),
    @(42602.577187499999, "Named", 2233, 408, 16, 3, 8, 27, 72, 0, 1, 0, 100)
)

$startRow = 30
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 1).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($r, 2).Value = $data[1]
    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
